$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (style) of the last existing header cell (G1) onto the new
# header cell (H1), then set its text.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "consequents_length"

# Fill the new "consequents_length" column (H2:H10) with the value 1 for
# every data row.
for ($r = 2; $r -le 10; $r++) {
    $ws.Cells.Item($r, 8).Value = 1
}
